$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they remain strings
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '25.768.01'
$ws.Range("E2").Value = '  +5.89%  '
$ws.Range("D3").Value = '1.707.88'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("D5").Value = '330.84'
$ws.Range("E5").Value = '  +6.58%  '
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.3686'
$ws.Range("E7").Value = '  +1.24%  '
$ws.Range("D8").Value = '48.44'
$ws.Range("E8").Value = '  +3.86%  '
$ws.Range("D9").Value = '0.3306'
$ws.Range("E9").Value = '  +2.58%  '
$ws.Range("D10").Value = '1.171'
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").Value = '0.07357'
$ws.Range("E11").Value = '  +5.44%  '
$ws.Range("D12").Value = '1.0000'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '6.205'
$ws.Range("E13").Value = '  +5.21%  '
$ws.Range("D14").Value = '19.98'
$ws.Range("E14").Value = '  +4.15%  '
$ws.Range("D15").Value = '6.870'
$ws.Range("E15").Value = '  +5.11%  '
$ws.Range("D16").Value = '1.699.05'
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").Value = '0.00001072'
$ws.Range("E17").Value = '  +4.48%  '
$ws.Range("D18").Value = '0.06618'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '81.40'
$ws.Range("E19").Value = '  +4.94%  '
$ws.Range("D20").Value = '0.9985'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '6.069'
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("E22").Value = '  +4.86%  '
$ws.Range("E23").Value = '  +4.75%  '
$ws.Range("D24").Value = '25.741.72'
$ws.Range("E24").Value = '  +5.74%  '
$ws.Range("D25").Value = '2.469'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '2.482'
$ws.Range("E26").Value = '  +9.34%  '
$ws.Range("D27").Value = '149.82'
$ws.Range("E27").Value = '  +3.07%  '
$ws.Range("D28").Value = '19.13'
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("D29").Value = '1.305'
$ws.Range("E29").Value = '  +12.76%  '
$ws.Range("D30").Value = '1.888.22'
$ws.Range("E30").Value = '  +3.35%  '
$ws.Range("D31").Value = '128.15'
$ws.Range("E31").Value = '  +4.02%  '
$ws.Range("D32").Value = '4.114'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").Value = '5.960'
$ws.Range("E33").Value = '  +6.81%  '
$ws.Range("D34").Value = '0.08481'
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '12.91'
$ws.Range("E35").Value = '  +8.33%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.677'
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").Value = '5.322'
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("D38").Value = '1.271'
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").Value = '0.06215'
$ws.Range("E39").Value = '  +4.80%  '
$ws.Range("D40").Value = '8.548'
$ws.Range("E40").Value = '  +6.61%  '
$ws.Range("D41").Value = '0.2122'
$ws.Range("E41").Value = '  +4.36%  '
$ws.Range("D42").Value = '0.02257'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("D43").Value = '14.49'
$ws.Range("E43").Value = '  +16.26%  '
$ws.Range("D44").Value = '0.6128'
$ws.Range("E44").Value = '  +5.25%  '
$ws.Range("D45").Value = '0.9989'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '3.847'
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").Value = '0.5841'
$ws.Range("E47").Value = '  +5.56%  '
$ws.Range("D48").Value = '126.40'
$ws.Range("E48").Value = '  +4.02%  '
$ws.Range("D49").Value = '2.005'
$ws.Range("E49").Value = '  +4.25%  '
$ws.Range("D50").Value = '0.07217'
$ws.Range("D51").Value = '1.207'
$ws.Range("E51").Value = '  +2.92%  '
